# Maestro.xlsx update: insert new "Agua Saborizada pera" (Levite) article
# as row 6, shifting all subsequent article rows down by one (rows 6-37 -> 7-38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents (in row order) for rows 6..38 after the insertion.
$data = @(
  ,@(7799155000579, 'Agua', 'Saborizada', 'pera', 'Levite', 500, 'ml.', 'Botella', 'Aguas saborizadas', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7799155000579.png', $true)
  ,@(1111, 'Salame', 'pelado', 'milán', 'Fela', 1000, 'gr.', 'Horma', 'Fiambres', 'Argentina', 1, $true, $true, 'C:\VentaSoft\Imágenes de artículos\1111.png', $true)
  ,@(2222, 'Salame', 'tipo', 'crespón', 'Paladini', 1000, 'gr.', 'Horma', 'Fiambres', 'Argentina', 1, $true, $true, 'C:\VentaSoft\Imágenes de artículos\2222.png', $true)
  ,@(3333, 'Queso', 'tipo', 'cremoso', 'La Paulina', 1000, 'gr.', 'Horma', 'Quesos', 'Argentina', 1, $true, $true, 'C:\VentaSoft\Imágenes de artículos\3333.png', $true)
  ,@(77990112, 'Licor', 'crema', 'café al coñac', 'Tres Plumas', 200, 'ml.', 'Petaca', 'Licores', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\77990112.png', $true)
  ,@(3086123206873, 'Encendedor', 'no recargable', 'colores varios', 'BX7', 1, 'und.', 'Suelto', 'Encendedores', 'Francia', 50, $false, $false, 'C:\VentaSoft\Imágenes de artículos\3086123206873.png', $true)
  ,@(7790036973036, 'Puré de tomates', 'sin conservantes', 'de la huerta', 'Baggio', 210, 'gr.', 'Tetra Brik', 'Puré de tomates', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790036973036.png', $true)
  ,@(7790070336385, 'Fideos', 'secos', 'spaghetti Nº 7', 'Luchetti', 500, 'gr.', 'Bolsa', 'Fideos secos', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790070336385.png', $true)
  ,@(7790150006153, 'Café', 'torrado molido', 'clásico', 'La Virginia', 500, 'gr.', 'Bolsa', 'Cafés', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790150006153.png', $true)
  ,@(7790150430392, 'Bicarbonato', 'de sodio', 'en bolsa', 'Alicante', 50, 'gr.', 'Bolsa', 'Especias', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790150430392.png', $true)
  ,@(7790520010445, 'Espirales insecticidas', 'repelentes', 'contra mosquitos', 'Raid', 12, 'und.', 'Caja', 'Insecticidas', 'Indonesia', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790520010445.png', $true)
  ,@(7790520016461, 'Desodorante', 'para ambientes en aerosol', 'caricias de algodón', 'Glade', 360, 'cm3.', 'Lata', 'Desodorantes', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790520016461.png', $true)
  ,@(7790520995360, 'Tabletas insecticidas', 'repelentes', 'contra mosquitos', 'Raid', 24, 'und.', 'Caja', 'Insecticidas', 'Argentina', 24, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790520995360.png', $true)
  ,@(7790639003574, 'Gaseosa', 'indian tonic', 'classic', 'Cunnington', 1.5, 'lt.', 'Botella', 'Gaseosas', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790639003574.png', $true)
  ,@(7791070000382, 'Papel higiénico', 'simple hoja', '"soft"', 'Ecco', 12, 'und.', 'Bolsa', 'Papeles Higiénicos', 'Argentina', 4, $false, $false, 'C:\VentaSoft\Imágenes de artículos\7791070000382.png', $false)
  ,@(7791600087128, 'Desodorante', 'para hombre', 'US', 'Colbert', 150, 'ml.', 'Lata', 'Desodorantes', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791600087128.png', $true)
  ,@(7791664000156, 'Ravioles', 'rellenos de', 'carne y espinaca', 'La Italiana', 500, 'gr.', 'Blister', 'Pastas frescas', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791664000156.png', $true)
  ,@(7791664000453, 'Ñoquis', 'frescos', 'de papa', 'La Italiana', 500, 'gr.', 'Blister', 'Pastas frescas', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791664000453.png', $true)
  ,@(7791684000934, 'Ravioles', 'rellenos de', 'jamón y muzzarella', 'La Italiana', 500, 'gr.', 'Blister', 'Pastas frescas', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791684000934.png', $true)
  ,@(7793913000139, 'Queso', 'rallado', 'sin gluten/TACC', 'Tregar', 1, 'und.', 'Sobre', 'Quesos', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7793913000139.png', $true)
  ,@(7795018002902, 'Leche', 'en polvo', 'entera', 'Quelech', 800, 'gr.', 'Bolsa', 'Papeles Higiénicos', 'Argentina', 12, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7795018002902.png', $true)
  ,@(7798100200491, 'Cigarrillos', 'rubios mentolados', 'en caja', 'Milenio', 20, 'und.', 'Caja', 'Cigarrillos', 'Argentina', 10, $false, $false, 'C:\VentaSoft\Imágenes de artículos\7798100200491.png', $true)
  ,@(7798100200583, 'Cigarrillos', 'rubios convertibles', 'en caja', 'Mill', 20, 'und.', 'Caja', 'Cigarrillos', 'Argentina', 10, $false, $false, 'C:\VentaSoft\Imágenes de artículos\7798100200583.png', $true)
  ,@(7791600087012, 'Eau de toilette', 'para hombre', 'US', 'Colbert', 60, 'ml.', 'Caja', 'Perfumes/Colonias', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791600087012.png', $true)
  ,@(7791520009743, 'Talco', 'para hombre', 'original', 'Veritas', 180, 'gr.', 'Botella', 'Desodorantes', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791520009743.png', $true)
  ,@(7794000006065, 'Mayonesa', 'libre de gluten', 'clásica', 'Hellmann''s', 237, 'gr.', 'Pouch', 'Aderezos', 'Argentina', 24, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7794000006065.png', $true)
  ,@(7792540260138, 'Azúcar', 'común tipo "A"', 'clásica', 'Ledesma', 1, 'kg.', 'Bolsa', 'Azúcar', 'Argentina', 10, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7792540260138.png', $false)
  ,@(7790639003895, 'Gaseosa', 'indian tonic', 'suave', 'Cunnington', 1.5, 'lt.', 'Botella', 'Gaseosas', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7790639003895.png', $true)
  ,@(7791600174767, 'Desodorante', 'para hombre', 'noir', 'Colbert', 250, 'ml.', 'Lata', 'Desodorantes', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791600174767.png', $true)
  ,@(7791600192372, 'Desodorante', 'para hombre', 'verde', 'Colbert', 150, 'ml.', 'Lata', 'Desodorantes', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791600192372.png', $true)
  ,@(7791600192488, 'Desodorante', 'para hombre', 'verde', 'Colbert', 250, 'ml.', 'Lata', 'Desodorantes', 'Argentina', 6, $false, $true, 'C:\VentaSoft\Imágenes de artículos\7791600192488.png', $true)
  ,@(2001300000798, 'Huevos', 'frescos medianos', 'de granja', 'Cismondi', 6, 'und.', 'Maple', 'Huevos', 'Argentina', 1, $false, $true, 'C:\VentaSoft\Imágenes de artículos\2001300000798.png', $false)
  ,@(7791070000078, 'Papel higiénico', 'simple hoja', '"soft"', 'Campanita', 4, 'und.', 'Bolsa', 'Papeles Higiénicos', 'Argentina', 10, $false, $false, 'C:\VentaSoft\Imágenes de artículos\7791070000078.png', $true)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $excelRow = 6 + $i
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
  }
}

# Row 38 is brand new (sheet previously ended at row 37); give column A the
# same integer number-format style ("0") already used by every other A-column cell.
$ws.Cells.Item(38, 1).NumberFormat = "0"
